$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Issue with 3.3v regulator overpowering:
# - rename the bottle cap part to clarify diameter
# - rename sonar sensor / LED parts and add a note cell for the ultrasonic sensor

$ws.Range("A3").Value = "1 inch diameter bottle cap"
$ws.Range("A17").Value = "LED Colored"
$ws.Range("A16").Value = "Sensor Sonar "
$ws.Range("F16").Value = "Sensor-Ultrasonic used for angry engineers"

$ws.Range("A1:E21").Select()
